# Update the metrics_11_2 results sheet: the 25 trained models ("model_11_2_0"
# .. "model_11_2_24") are re-ordered in column A (best models first), and
# every data row B2:I26 is refreshed with the newly-computed metric values
# for the current best model (r2, r2_test, r2_val, r2_vt, mse, mse_test,
# mse_val, mse_vt) - matching the commit "4 mdelo melhores rstds".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of model names for rows 2..26 (A column)
$modelNames = @(
    "model_11_2_0",
    "model_11_2_22",
    "model_11_2_21",
    "model_11_2_20",
    "model_11_2_19",
    "model_11_2_18",
    "model_11_2_17",
    "model_11_2_16",
    "model_11_2_15",
    "model_11_2_14",
    "model_11_2_13",
    "model_11_2_23",
    "model_11_2_12",
    "model_11_2_10",
    "model_11_2_9",
    "model_11_2_8",
    "model_11_2_7",
    "model_11_2_6",
    "model_11_2_5",
    "model_11_2_4",
    "model_11_2_3",
    "model_11_2_2",
    "model_11_2_1",
    "model_11_2_11",
    "model_11_2_24"
)

# Updated metric values shared by every row: r2, r2_test, r2_val, r2_vt, mse,
# mse_test, mse_val, mse_vt (columns B..I)
$metricValues = @(
    0.3494677884409869,
    0.267500787098648,
    0.3229168344848683,
    0.4396081365611429,
    0.7199474573135376,
    0.8601891994476318,
    0.5655463933944702,
    0.7215338945388794
)

$firstRow = 2
for ($i = 0; $i -lt $modelNames.Length; $i++) {
    $row = $firstRow + $i

    $ws.Cells.Item($row, 1).Value = $modelNames[$i]

    for ($c = 0; $c -lt $metricValues.Length; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $metricValues[$c]
    }
}
